$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column cells to retain plain-text storage (matching the
# original inline-string cells) instead of Excel auto-coercing numeric-looking
# text into real numbers. We set an explicit "@" (Text) format before writing,
# then reset the style back to "Normal" so the cells end up with no style index
# attribute at all -- exactly like the source file.

$priceCells = @("D2","D3","D5","D10","D11","D12","D14","D15","D16","D17","D19","D21","D22","D24","D25","D27","D28","D29","D31","D35","D38","D41","D44","D45","D46","D48","D49","D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.762.50'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '1.650.33'
$ws.Range("E3").Value = '  +0.88%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '215.91'
$ws.Range("E5").Value = '  +1.36%  '
$ws.Range("E6").Value = '  +0.95%  '
$ws.Range("E8").Value = '  -0.45%  '
$ws.Range("E9").Value = '  +0.99%  '
$ws.Range("D10").Value = '19.40'
$ws.Range("E10").Value = '  +1.51%  '
$ws.Range("D11").Value = '0.0846'
$ws.Range("E11").Value = '  +0.49%  '
$ws.Range("D12").Value = '1.879.95'
$ws.Range("E12").Value = '  +0.81%  '
$ws.Range("E13").Value = '  +2.94%  '
$ws.Range("D14").Value = '1.639.04'
$ws.Range("E14").Value = '  +0.60%  '
$ws.Range("D15").Value = '0.535'
$ws.Range("E15").Value = '  +1.42%  '
$ws.Range("D16").Value = '66.62'
$ws.Range("E16").Value = '  +5.23%  '
$ws.Range("D17").Value = '26.800.87'
$ws.Range("E17").Value = '  +0.51%  '
$ws.Range("E18").Value = '  +1.69%  '
$ws.Range("D19").Value = '221.04'
$ws.Range("E19").Value = '  +1.56%  '
$ws.Range("E20").Value = '  +0.06%  '
$ws.Range("D21").Value = '4.41'
$ws.Range("E21").Value = '  +2.29%  '
$ws.Range("D22").Value = '6.39'
$ws.Range("E22").Value = '  +2.64%  '
$ws.Range("E23").Value = '  +0.71%  '
$ws.Range("D24").Value = '2.16'
$ws.Range("E24").Value = '  +12.52%  '
$ws.Range("D25").Value = '147.49'
$ws.Range("E25").Value = '  -1.18%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("D27").Value = '0.122'
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("D28").Value = '7.10'
$ws.Range("E28").Value = '  +3.18%  '
$ws.Range("D29").Value = '15.98'
$ws.Range("E29").Value = '  +3.36%  '
$ws.Range("E30").Value = '  +0.50%  '
$ws.Range("D31").Value = '1.18'
$ws.Range("E31").Value = '  +0.64%  '
$ws.Range("E32").Value = '  +4.71%  '
$ws.Range("E33").Value = '  +4.15%  '
$ws.Range("E34").Value = '  +4.71%  '
$ws.Range("D35").Value = '1.299.47'
$ws.Range("E35").Value = '  +9.17%  '
$ws.Range("E36").Value = '  +5.32%  '
$ws.Range("E37").Value = '  +0.91%  '
$ws.Range("D38").Value = '0.834'
$ws.Range("E38").Value = '  +2.82%  '
$ws.Range("E39").Value = '  +3.93%  '
$ws.Range("D41").Value = '0.817'
$ws.Range("E41").Value = '  +2.85%  '
$ws.Range("E42").Value = '  -2.59%  '
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("D44").Value = '1.792.22'
$ws.Range("E44").Value = '  +1.08%  '
$ws.Range("D45").Value = '94.01'
$ws.Range("E45").Value = '  +1.52%  '
$ws.Range("D46").Value = '61.13'
$ws.Range("E46").Value = '  +11.37%  '
$ws.Range("E47").Value = '  +4.84%  '
$ws.Range("D48").Value = '0.0518'
$ws.Range("E48").Value = '  +0.89%  '
$ws.Range("D49").Value = '7.84'
$ws.Range("E49").Value = '  +2.65%  '
$ws.Range("D50").Value = '0.0981'
$ws.Range("E50").Value = '  +3.62%  '
$ws.Range("E51").Value = '  -0.61%  '

# Reset style on the touched Price cells back to Normal so no stray style index
# is left on the cell (the NumberFormat="@" above otherwise sticks as "s=...").
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
